# Apply the "corrected ICDC Breed 1-14 scripts" edit.
#
# The file-level Cypher query stored in cell B4 of the "startup" sheet is
# corrected: two projected columns (`f.file_type` AS `File Type` and
# `demo.breed` AS Breed) are removed from the RETURN clause, and the
# wrapped-text row shrinks accordingly (two fewer lines of text).
# B2 (the case-level query) is left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newFileQuery = "MATCH (f:file)-->(parent)`n" + `
    "WITH DISTINCT f, parent`n" + `
    "MATCH (f)-[*]->(c:case)<--(demo:demographic)`n" + `
    "WHERE demo.breed IN ['Chihuahua']  `n" + `
    "OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n" + `
    "OPTIONAL MATCH (samp:sample)-->(c)`n" + `
    "WITH DISTINCT f, parent, c, demo, diag, s`n" + `
    "RETURN  coalesce(f.file_name, '') AS ``File Name``,`n" + `
    "        coalesce(labels(parent)[0], '') AS ``Association``,`n" + `
    "        coalesce(f.file_description, '') AS ``Description``,`n" + `
    "        coalesce(f.file_format, '') AS ``Format``,`n" + `
    "        coalesce(f.file_size, '') AS ``Size``,`n" + `
    "        coalesce(c.case_id, '') AS ``Case ID``,`n" + `
    "        coalesce(diag.disease_term,'') AS Diagnosis , `n" + `
    "        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Cells.Item(4, 2).Value2 = $newFileQuery

# The row shrinks by two wrapped lines now that the query text is shorter.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection / scroll moves down onto the corrected cell.
$ws.Activate()
$ws.Range("B4").Select()
